$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the existing row 143, shifting the old rows
# 143-158 down to 147-162 (their content is unchanged by this edit).
$ws.Rows("143:146").Insert()

# New row 143: Sandia, Cuarta, Región de O'Higgins, week of 44551
$ws.Range("A143").Value = 4
$ws.Range("B143").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C143").Value = 'Los Lagos'
$ws.Range("D143").Value = 44551
$ws.Range("E143").Value = 10
$ws.Range("F143").Value = 100112028
$ws.Range("G143").Value = 'Sandia'
$ws.Range("H143").Value = 'Sin especificar'
$ws.Range("I143").Value = 'Cuarta'
$ws.Range("J143").Value = 600
$ws.Range("K143").Value = 1500
$ws.Range("L143").Value = 1500
$ws.Range("M143").Value = 1500
$ws.Range("N143").Value = '$/unidad'
$ws.Range("O143").Value = 'Región de O''Higgins'
$ws.Range("P143").Value = 1500
$ws.Range("Q143").Value = 1
$ws.Range("R143").Value = 'Hortaliza'

# New row 144: Sandia, Primera, Región de O'Higgins, week of 44551
$ws.Range("A144").Value = 4
$ws.Range("B144").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C144").Value = 'Los Lagos'
$ws.Range("D144").Value = 44551
$ws.Range("E144").Value = 10
$ws.Range("F144").Value = 100112028
$ws.Range("G144").Value = 'Sandia'
$ws.Range("H144").Value = 'Sin especificar'
$ws.Range("I144").Value = 'Primera'
$ws.Range("J144").Value = 400
$ws.Range("K144").Value = 3000
$ws.Range("L144").Value = 3000
$ws.Range("M144").Value = 3000
$ws.Range("N144").Value = '$/unidad'
$ws.Range("O144").Value = 'Región de O''Higgins'
$ws.Range("P144").Value = 3000
$ws.Range("Q144").Value = 1
$ws.Range("R144").Value = 'Hortaliza'

# New row 145: Sandia, Segunda, Región de O'Higgins, week of 44551
$ws.Range("A145").Value = 4
$ws.Range("B145").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C145").Value = 'Los Lagos'
$ws.Range("D145").Value = 44551
$ws.Range("E145").Value = 10
$ws.Range("F145").Value = 100112028
$ws.Range("G145").Value = 'Sandia'
$ws.Range("H145").Value = 'Sin especificar'
$ws.Range("I145").Value = 'Segunda'
$ws.Range("J145").Value = 400
$ws.Range("K145").Value = 2500
$ws.Range("L145").Value = 2500
$ws.Range("M145").Value = 2500
$ws.Range("N145").Value = '$/unidad'
$ws.Range("O145").Value = 'Región de O''Higgins'
$ws.Range("P145").Value = 2500
$ws.Range("Q145").Value = 1
$ws.Range("R145").Value = 'Hortaliza'

# New row 146: Sandia, Tercera, Región de O'Higgins, week of 44551
$ws.Range("A146").Value = 4
$ws.Range("B146").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C146").Value = 'Los Lagos'
$ws.Range("D146").Value = 44551
$ws.Range("E146").Value = 10
$ws.Range("F146").Value = 100112028
$ws.Range("G146").Value = 'Sandia'
$ws.Range("H146").Value = 'Sin especificar'
$ws.Range("I146").Value = 'Tercera'
$ws.Range("J146").Value = 400
$ws.Range("K146").Value = 2000
$ws.Range("L146").Value = 2000
$ws.Range("M146").Value = 2000
$ws.Range("N146").Value = '$/unidad'
$ws.Range("O146").Value = 'Región de O''Higgins'
$ws.Range("P146").Value = 2000
$ws.Range("Q146").Value = 1
$ws.Range("R146").Value = 'Hortaliza'
